# =====================================================================
# Edit: "modificación en nombres de columnas, agregada columna de
#        relación con obras en referentes"
#
# Summary of what this script does:
#  1. Renames sheet "obras" -> "obra".
#  2. Lower-cases / renames the header row of both sheets
#     (Creator/Title/Date/ancho cm/alto cm/Format.medium/Archivo/Referentes/
#      Publisher -> creator/title/date/ancho/alto/medium/archivo/
#      referentes/publisher), keeping "ID" as-is.
#  3. Converts the plain "year" numbers in "obra"!D2:D4 into real dates
#     (1965, 1997, 1983 -> Jan 1 of that year) formatted as yyyy-mm-dd,
#     matching the date format already used in "refentes" (now also
#     yyyy-mm-dd instead of the old default date format).
#  4. Adds a new "obra" column (F) to "refentes" that links each
#     "referente" row back to the id of the related row in "obra".
#  5. Leaves the selection / active sheet the way the author left it
#     (refentes active, D4 selected in "obra", B5 selected in "refentes").
# =====================================================================

$wb = $excel.ActiveWorkbook

$wsObra = $wb.Worksheets.Item(1)
$wsRef  = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# 1. Rename the first sheet
# ---------------------------------------------------------------
$wsObra.Name = "obra"

# ---------------------------------------------------------------
# 2. Header renames - sheet "obra"
#    A:id  B:creator  C:title  D:date  E:ancho  F:alto  G:medium  H:archivo  I:referentes
#    (written in the same order the author retyped them in, D..I first,
#    then B, C)
# ---------------------------------------------------------------
$wsObra.Range("A1").Value2 = "ID"
$wsObra.Range("D1").Value2 = "date"
$wsObra.Range("E1").Value2 = "ancho"
$wsObra.Range("F1").Value2 = "alto"
$wsObra.Range("G1").Value2 = "medium"
$wsObra.Range("H1").Value2 = "archivo"
$wsObra.Range("I1").Value2 = "referentes"
$wsObra.Range("B1").Value2 = "creator"
$wsObra.Range("C1").Value2 = "title"

# Header row keeps its bold style; D1 (date header) additionally gets the
# custom date format applied (bold + yyyy-mm-dd), matching "refentes"!C1.
$wsObra.Range("D1").NumberFormat = "yyyy\-mm\-dd;@"

# ---------------------------------------------------------------
# 3. Convert the "date" column (D2:D4) of "obra" from bare years to
#    real date values, formatted yyyy-mm-dd.
# ---------------------------------------------------------------
$wsObra.Range("D2:D4").NumberFormat = "yyyy\-mm\-dd;@"
$wsObra.Range("D2").Value2 = 23743   # 1965-01-01
$wsObra.Range("D3").Value2 = 35431   # 1997-01-01
$wsObra.Range("D4").Value2 = 30317   # 1983-01-01

# ---------------------------------------------------------------
# 4. Header renames - sheet "refentes"
#    A:id  B:title  C:date  D:publisher  E:archivo  F:obra(new)
# ---------------------------------------------------------------
$wsRef.Range("A1").Value2 = "ID"
$wsRef.Range("C1").Value2 = "date"
$wsRef.Range("B1").Value2 = "title"
$wsRef.Range("D1").Value2 = "publisher"
$wsRef.Range("E1").Value2 = "archivo"

$wsRef.Range("C1").NumberFormat = "yyyy\-mm\-dd;@"
$wsRef.Range("C2:C3").NumberFormat = "yyyy\-mm\-dd;@"

# New column F "obra": links each referente back to the related obra id.
$wsRef.Range("F1").Value2 = "obra"
$wsRef.Range("F2").Value2 = 1
$wsRef.Range("F3").Value2 = 2
$wsRef.Range("F4").Value2 = 2
$wsRef.Range("F5").Value2 = 3

# ---------------------------------------------------------------
# 5. Selection / active sheet state
# ---------------------------------------------------------------
$wsObra.Range("D4").Select() | Out-Null
$wsRef.Range("B5").Select() | Out-Null
$wsRef.Select() | Out-Null
